$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 200-201, shifting existing rows 200-248 down to 202-250.
$ws.Rows("200:201").Insert()

# --- New row 200 ---
$ws.Range("A200").Value = 3
$ws.Range("B200").Value = "Femacal de La Calera"
$ws.Range("C200").Value = "Coquimbo"
$ws.Range("D200").Value = 44511
$ws.Range("E200").Value = 5
$ws.Range("F200").Value = 100112013
$ws.Range("G200").Value = "Alcachofa"
$ws.Range("H200").Value = "Española"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 13000
$ws.Range("K200").Value = 300
$ws.Range("L200").Value = 320
$ws.Range("M200").Value = 310
$ws.Range("N200").Value = "$/unidad"
$ws.Range("O200").Value = "Llay Llay"
$ws.Range("P200").Value = 310
$ws.Range("Q200").Value = 1
$ws.Range("R200").Value = "Hortaliza"

# --- New row 201 ---
$ws.Range("A201").Value = 3
$ws.Range("B201").Value = "Femacal de La Calera"
$ws.Range("C201").Value = "Coquimbo"
$ws.Range("D201").Value = 44511
$ws.Range("E201").Value = 5
$ws.Range("F201").Value = 100112013
$ws.Range("G201").Value = "Alcachofa"
$ws.Range("H201").Value = "Española"
$ws.Range("I201").Value = "Segunda"
$ws.Range("J201").Value = 6800
$ws.Range("K201").Value = 200
$ws.Range("L201").Value = 200
$ws.Range("M201").Value = 200
$ws.Range("N201").Value = "$/unidad"
$ws.Range("O201").Value = "Llay Llay"
$ws.Range("P201").Value = 200
$ws.Range("Q201").Value = 1
$ws.Range("R201").Value = "Hortaliza"
